# Adds a new results row (row 25) to Sheet1 for the
# "Bandpower + PCA + NuSVM (poly kernel)" method, matching the author's
# commit ("Plotting some graphs and figures").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of results, appended right after the existing last row (24).
$ws.Cells.Item(25, 1).Value = "Bandpower + PCA + NuSVM (poly kernel)"

$ws.Cells.Item(25, 2).Value = 0.9456
$ws.Cells.Item(25, 2).NumberFormat = "0.00%"

$ws.Cells.Item(25, 3).Value = "19/19"
$ws.Cells.Item(25, 4).Value = "RH"
$ws.Cells.Item(25, 5).Value = "3, 5, 11, 13"
$ws.Cells.Item(25, 6).Value = "nu=865, n_components=3, freq bands (Hz) 4-8,8-13,13-30"

# Reflect the scrolled/selected view from the source edit (best effort -
# harmless if the host doesn't persist window scroll position).
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 3 } catch {}

$ws.Range("F25").Select()
